$d = $word.ActiveDocument

function Split-AtBoundary($findText, $replaceText) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# --- "Programa" paragraph: split the 7 numbered items with manual line breaks ---
Split-AtBoundary "elabora$([char]0x00E7)$([char]0x00E3)o de projeto;2.Metodologia" "elabora$([char]0x00E7)$([char]0x00E3)o de projeto;^l2.Metodologia"
Split-AtBoundary "teste do produto;3. Processo" "teste do produto;^l3. Processo"
Split-AtBoundary "Checar-Agir;4.M$([char]0x00E9)todos" "Checar-Agir;^l4.M$([char]0x00E9)todos"
Split-AtBoundary "cient$([char]0x00ED)ficos;5.Desenvolvimento" "cient$([char]0x00ED)ficos;^l5.Desenvolvimento"
Split-AtBoundary "especifica$([char]0x00E7)$([char]0x00E3)o da solu$([char]0x00E7)$([char]0x00E3)o;6.No$([char]0x00E7)$([char]0x00F5)es" "especifica$([char]0x00E7)$([char]0x00E3)o da solu$([char]0x00E7)$([char]0x00E3)o;^l6.No$([char]0x00E7)$([char]0x00F5)es"
Split-AtBoundary "equipes e times7.Tutoria" "equipes e times^l7.Tutoria"

# --- "M$([char]0x00E9)todo:" run: split the 4 sentences with manual line breaks ---
Split-AtBoundary "dentre outros.Os alunos" "dentre outros.^lOs alunos"
Split-AtBoundary "sua profiss$([char]0x00E3)o.Cada grupo" "sua profiss$([char]0x00E3)o.^lCada grupo"
Split-AtBoundary "do projeto.As aulas ocorrer$([char]0x00E3)o" "do projeto.^lAs aulas ocorrer$([char]0x00E3)o"

# --- "Crit$([char]0x00E9)rio:" run: split the 2 sentences with a manual line break ---
Split-AtBoundary "dentre outros.O detalhamento" "dentre outros.^lO detalhamento"

# --- "Bibliografia" paragraph: split the 5 reference entries with manual line breaks ---
Split-AtBoundary "2013.- COCIAN" "2013.^l- COCIAN"
Split-AtBoundary "2017.- BENNETT" "2017.^l- BENNETT"
Split-AtBoundary "2014.- ALEXANDER" "2014.^l- ALEXANDER"
Split-AtBoundary "2015.- MCCAHAN" "2015.^l- MCCAHAN"

Write-Host "All splits applied"
